# "added effort estimates to historical work"
#
# - Renames the original sheet "Sheet1" -> "Estimates"
# - Adds a new "Effort" worksheet (right after Estimates) with a header
#   row (bold) + one data row of effort estimates
# - Leaves the Estimates sheet's selection at C33 (no longer the
#   tab-selected sheet) and makes Effort the active/selected sheet with
#   column G selected

$wb = $excel.ActiveWorkbook

# --- rename existing sheet, add the new one right after it -----------
$estimates = $wb.Worksheets.Item(1)
$estimates.Name = "Estimates"

$effort = $wb.Worksheets.Add($null, $estimates)
$effort.Name = "Effort"

# --- header row (bold) --------------------------------------------------
$headers = @("year", "nvessels", "mesh_in_avg", "sets_tot_est", "sets_obs", "obs_perc", "notes")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $effort.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
}

# --- data row ------------------------------------------------------------
$effort.Cells.Item(2, 1).Value = 2010
$effort.Cells.Item(2, 2).Value = 50
$effort.Cells.Item(2, 3).Value = 7.2
$effort.Cells.Item(2, 4).Value = 1724
$effort.Cells.Item(2, 5).Value = 216
$effort.Cells.Item(2, 6).Value = 12.5
$effort.Cells.Item(2, 7).Value = "ests are from 2009 logbook data"

# widen the notes column to fit its text
$effort.Columns.Item(7).ColumnWidth = 27.33

# --- view/selection state -------------------------------------------------
$estimates.Range("C33").Select()
$effort.Activate()
$effort.Columns.Item(7).Select()
